# "now works with familiarisation area"
# Update the assignment test sheet so the assignment/partner table also
# covers the "familiarisation area" placeholder students (FA_Student2..6),
# refresh the open/close dates, widen the table to include the two new
# (currently blank) columns E and F, and re-fit the columns to the new
# (longer) content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: refresh the open / close dates (kept as real date serials) ---
$ws.Range("C1").Value = 45649
$ws.Range("D1").Value = 45656

# --- Column A: swap the numeric student ids for the familiarisation-area
# placeholder names. Enter A3 before A2 so the new shared strings land in
# the same order as the source file (FA_Student2 first, then FA_Student3).
$ws.Range("A3").Value = "FA_Student2"
$ws.Range("A2").Value = "FA_Student3"
$ws.Range("A4").Value = "FA_Student4"
$ws.Range("A5").Value = "FA_Student5"
$ws.Range("A6").Value = "FA_Student6"

# Row 6 no longer wraps onto two lines now that it holds "FA_Student6"
# instead of the long numeric id, so its row shrinks back to single-line
# height just like the rows above it.
$ws.Rows(6).RowHeight = 16.5

# A highlight (white fill) was left on the familiarisation row (row 3).
$ws.Range("A3").Interior.Color = 16777215

# The other familiarisation rows (2, 4, 5) lost their right/center/wrap
# alignment on column A - copy the neighbouring cell's font formatting in,
# then drop the alignment back to the default.
foreach ($r in 2, 4, 5) {
    $ws.Range("B$r").Copy() | Out-Null
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("A$r").WrapText = $false
    $ws.Range("A$r").HorizontalAlignment = 1
    $ws.Range("A$r").VerticalAlignment = -4107
}

# --- New columns E and F: blank cells that pick up the same per-row
# formatting as column A in each row.
foreach ($r in 2, 3, 4, 5, 6) {
    $ws.Range("A$r").Copy() | Out-Null
    $ws.Range("E$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("F$r").PasteSpecial(-4122) | Out-Null
}
$ws.Range("E2:F2,E4:F5").WrapText = $false
$ws.Range("E2:F2,E4:F5").HorizontalAlignment = 1
$ws.Range("E2:F2,E4:F5").VerticalAlignment = -4107

# --- Re-fit the columns now that the table holds the new (differently
# sized) content.
$ws.Range("A1:F6").Columns.AutoFit()

# Leave the cursor where the author left it.
$ws.Range("D2").Select()
